# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the zh-cn and
# de-de report rows that correspond to the
# f9b19333-f2c6-44ca-be66-6ef7624ed513 handoff/handback pair, to reflect
# the freshly (re-)generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 04:42:59"
$wsZhCn.Range("H2").Value = "2016-03-14 04:43:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 04:43:02"
$wsDeDe.Range("H2").Value = "2016-03-14 04:43:22"
